$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the paragraph that begins "The results
#    of the three tests run with long.txt are as follows:" to the end of the
#    paragraph that ends in "...fast_stack and stackprof." (right after the
#    final period, before the paragraph mark).  Also drop the stray
#    <w:lastRenderedPageBreak/> that sat in front of the "long.txt" run.
# ---------------------------------------------------------------------------

# 1a. Delete the existing "_GoBack" bookmark (wherever Word currently has it).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 1b. Rebuild the "long.txt" run without the <w:lastRenderedPageBreak/> child,
#     keeping its text and run formatting intact.
$rLong = $d.Content
$rLong.Find.Execute("The results of the three tests run with long.txt are as follows:")
$rLongFixed = $d.Range($rLong.Start, $rLong.End)
$longXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>The results of the three tests run with long.txt are as follows:</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$rLongFixed.InsertXML($longXml)

# 1c. Re-add "_GoBack" right after the "." that ends "...stackprof." -
#     collapsed ranges that land exactly before a paragraph mark can't be
#     bookmarked directly, so nudge in a placeholder character, bookmark the
#     (non-collapsed) range around it, then delete the placeholder again.
$rStack = $d.Content
$rStack.Find.Execute("fast_stack and stackprof.")
$pos = $rStack.End
$rPlaceholder = $d.Range($pos, $pos)
$rPlaceholder.InsertAfter("X")
$rBookmark = $d.Range($pos, $pos + 1)
$d.Bookmarks.Add("_GoBack", $rBookmark)
$rCleanup = $d.Range($pos, $pos + 1)
$rCleanup.Text = ""

# ---------------------------------------------------------------------------
# 2) Replace the heading "#Edge Cases and Failure Modes" with the expanded
#    paragraph discussing the ArgumentCheck_test unit test edge cases.
# ---------------------------------------------------------------------------
$rHeading = $d.Content
$rHeading.Find.Execute("#Edge Cases and Failure Modes")
$rHeadingFixed = $d.Range($rHeading.Start, $rHeading.End)
$headingXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">In our unit tests, we checked for several edge cases. For example, in </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>ArgumentCheck_test</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> we checked for three edge cases in the number of arguments. The assignment required the program to only accept one argument, so we tested cases of 0,1, and 2 arguments. 0 and 2 were also failure cases. </w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$rHeadingFixed.InsertXML($headingXml)
